# Customize onClick function with Explictwait
#
# The "status"/"Pass" helper columns (E..U) get extended three more times
# to the right (V, W, X), mirroring the existing pattern: a "status"
# header cell (shared string) in row 1 with its own fill style, and
# "Pass" (shared string) in rows 2-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns right after the last "status"/"Pass" column (U)
# so the existing columns keep their values/styles untouched.
$ws.Columns("V:X").Insert()

# Same text used by every other status column.
$ws.Range("V1:X1").Value = "status"
$ws.Range("V2:X6").Value = "Pass"

# Give each new header cell its own fill (matching the look of the other
# status columns) so every column gets a distinct style entry, same as
# the rest of the sheet.
$ws.Range("V1").Interior.ColorIndex = 17
$ws.Range("W1").Interior.PatternColorIndex = 17
$ws.Range("X1").Interior.Color = 10079487

# Match the column width used by the rest of the status columns.
$ws.Columns("V:X").ColumnWidth = 5.5
